$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 137.0025913333334
$ws.Range("H2").Value = 411.007774
$ws.Range("I2").Value = 0.07043159922291199
$ws.Range("J2").Value = 0.07043159922291199
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 262.8951401926118
$ws.Range("R2").Value = 2366.056261733506
$ws.Range("S2").Value = 0.0004594971029764194
$ws.Range("T2").Value = 0.0004594971029764194
$ws.Range("G3").Value = 137.0025913333334
$ws.Range("H3").Value = 411.007774
$ws.Range("I3").Value = 0.07043159922291199
$ws.Range("J3").Value = 0.07043159922291199
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 24836.97939131808
$ws.Range("R3").Value = 223532.8145218627
$ws.Range("S3").Value = 0.04341092067595557
$ws.Range("T3").Value = 0.04341092067595557
$ws.Range("G4").Value = 137.0025913333334
$ws.Range("H4").Value = 411.007774
$ws.Range("I4").Value = 0.07043159922291199
$ws.Range("J4").Value = 0.07043159922291199
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 3979.606381868241
$ws.Range("R4").Value = 35816.45743681416
$ws.Range("S4").Value = 0.006955691923841492
$ws.Range("T4").Value = 0.006955691923841493
$ws.Range("G5").Value = 137.0025913333334
$ws.Range("H5").Value = 411.007774
$ws.Range("I5").Value = 0.07043159922291199
$ws.Range("J5").Value = 0.07043159922291199
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 11217.01939480149
$ws.Range("R5").Value = 100953.1745532134
$ws.Range("S5").Value = 0.01960548952013851
$ws.Range("T5").Value = 0.01960548952013851
$ws.Range("I6").Value = 0.1284841594777439
$ws.Range("J6").Value = 0.1284841594777439
$ws.Range("M6").Value = 1.918906333333333
$ws.Range("N6").Value = 5.756718999999999
$ws.Range("O6").Value = 0.006524019162508824
$ws.Range("P6").Value = 0.006524019162508824
$ws.Range("Q6").Value = 479.583901134012
$ws.Range("R6").Value = 4316.255110206108
$ws.Range("S6").Value = 0.000838233118511641
$ws.Range("T6").Value = 0.000838233118511641
$ws.Range("I7").Value = 0.1284841594777439
$ws.Range("J7").Value = 0.1284841594777439
$ws.Range("O7").Value = 0.6163557430885885
$ws.Range("P7").Value = 0.6163557430885885
$ws.Range("S7").Value = 0.07919194959001756
$ws.Range("T7").Value = 0.07919194959001756
$ws.Range("I8").Value = 0.1284841594777439
$ws.Range("J8").Value = 0.1284841594777439
$ws.Range("M8").Value = 29.04767233333333
$ws.Range("N8").Value = 87.143017
$ws.Range("O8").Value = 0.09875811426384234
$ws.Range("P8").Value = 0.09875811426384236
$ws.Range("Q8").Value = 7259.758214609316
$ws.Range("R8").Value = 65337.82393148384
$ws.Range("S8").Value = 0.01268885330279678
$ws.Range("T8").Value = 0.01268885330279678
$ws.Range("I9").Value = 0.1284841594777439
$ws.Range("J9").Value = 0.1284841594777439
$ws.Range("M9").Value = 81.87450533333333
$ws.Range("N9").Value = 245.623516
$ws.Range("O9").Value = 0.2783621234850603
$ws.Range("P9").Value = 0.2783621234850603
$ws.Range("Q9").Value = 20462.53847261477
$ws.Range("R9").Value = 184162.8462535329
$ws.Range("S9").Value = 0.03576512346641793
$ws.Range("T9").Value = 0.03576512346641793
$ws.Range("G10").Value = 88.73577866666666
$ws.Range("H10").Value = 266.207336
$ws.Range("I10").Value = 0.04561813567874526
$ws.Range("J10").Value = 0.04561813567874527
$ws.Range("M10").Value = 1.918906333333333
$ws.Range("N10").Value = 5.756718999999999
$ws.Range("O10").Value = 0.006524019162508824
$ws.Range("P10").Value = 0.006524019162508824
$ws.Range("Q10").Value = 170.2756476767315
$ws.Range("R10").Value = 1532.480829090584
$ws.Range("S10").Value = 0.0002976135913260615
$ws.Range("T10").Value = 0.0002976135913260616
$ws.Range("G11").Value = 88.73577866666666
$ws.Range("H11").Value = 266.207336
$ws.Range("I11").Value = 0.04561813567874526
$ws.Range("J11").Value = 0.04561813567874527
$ws.Range("O11").Value = 0.6163557430885885
$ws.Range("P11").Value = 0.6163557430885885
$ws.Range("Q11").Value = 16086.76656819072
$ws.Range("R11").Value = 144780.8991137165
$ws.Range("S11").Value = 0.02811699991458909
$ws.Range("T11").Value = 0.02811699991458909
$ws.Range("G12").Value = 88.73577866666666
$ws.Range("H12").Value = 266.207336
$ws.Range("I12").Value = 0.04561813567874526
$ws.Range("J12").Value = 0.04561813567874527
$ws.Range("M12").Value = 29.04767233333333
$ws.Range("N12").Value = 87.143017
$ws.Range("O12").Value = 0.09875811426384234
$ws.Range("P12").Value = 0.09875811426384236
$ws.Range("Q12").Value = 2577.567822952523
$ws.Range("R12").Value = 23198.11040657271
$ws.Range("S12").Value = 0.004505161055864988
$ws.Range("T12").Value = 0.004505161055864989
$ws.Range("G13").Value = 88.73577866666666
$ws.Range("H13").Value = 266.207336
$ws.Range("I13").Value = 0.04561813567874526
$ws.Range("J13").Value = 0.04561813567874527
$ws.Range("M13").Value = 81.87450533333333
$ws.Range("N13").Value = 245.623516
$ws.Range("O13").Value = 0.2783621234850603
$ws.Range("P13").Value = 0.2783621234850603
$ws.Range("Q13").Value = 7265.197983701486
$ws.Range("R13").Value = 65386.78185331338
$ws.Range("S13").Value = 0.01269836111696512
$ws.Range("T13").Value = 0.01269836111696512
$ws.Range("G14").Value = 1469.52242
$ws.Range("H14").Value = 4408.56726
$ws.Range("I14").Value = 0.7554661056205989
$ws.Range("J14").Value = 0.7554661056205988
$ws.Range("M14").Value = 1.918906333333333
$ws.Range("N14").Value = 5.756718999999999
$ws.Range("O14").Value = 0.006524019162508824
$ws.Range("P14").Value = 0.006524019162508824
$ws.Range("Q14").Value = 2819.875878713327
$ws.Range("R14").Value = 25378.88290841994
$ws.Range("S14").Value = 0.004928675349694702
$ws.Range("T14").Value = 0.004928675349694702
$ws.Range("G15").Value = 1469.52242
$ws.Range("H15").Value = 4408.56726
$ws.Range("I15").Value = 0.7554661056205989
$ws.Range("J15").Value = 0.7554661056205988
$ws.Range("O15").Value = 0.6163557430885885
$ws.Range("P15").Value = 0.6163557430885885
$ws.Range("Q15").Value = 266407.355550067
$ws.Range("R15").Value = 2397666.199950603
$ws.Range("S15").Value = 0.4656358729080263
$ws.Range("T15").Value = 0.4656358729080263
$ws.Range("G16").Value = 1469.52242
$ws.Range("H16").Value = 4408.56726
$ws.Range("I16").Value = 0.7554661056205989
$ws.Range("J16").Value = 0.7554661056205988
$ws.Range("M16").Value = 29.04767233333333
$ws.Range("N16").Value = 87.143017
$ws.Range("O16").Value = 0.09875811426384234
$ws.Range("P16").Value = 0.09875811426384236
$ws.Range("Q16").Value = 42686.20574264705
$ws.Range("R16").Value = 384175.8516838234
$ws.Range("S16").Value = 0.07460840798133909
$ws.Range("T16").Value = 0.07460840798133909
$ws.Range("G17").Value = 1469.52242
$ws.Range("H17").Value = 4408.56726
$ws.Range("I17").Value = 0.7554661056205989
$ws.Range("J17").Value = 0.7554661056205988
$ws.Range("M17").Value = 81.87450533333333
$ws.Range("N17").Value = 245.623516
$ws.Range("O17").Value = 0.2783621234850603
$ws.Range("P17").Value = 0.2783621234850603
$ws.Range("Q17").Value = 120316.4212137429
$ws.Range("R17").Value = 1082847.790923686
$ws.Range("S17").Value = 0.2102931493815388
$ws.Range("T17").Value = 0.2102931493815387
